# Insert a new weekly record at row 536 of the daily-logic Brócoli
# (Agrícola del Norte S.A. de Arica) consolidated sheet, pushing the
# existing rows 536-573 down to 537-574.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 536:573 down by one to make room for the new record.
$ws.Rows("536:536").Insert()

# Populate the newly inserted row 536 with the new observation.
$ws.Range("A536").Value = 1
$ws.Range("B536").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C536").Value = "Arica y Parinacota"
$ws.Range("D536").Value = 45106
$ws.Range("E536").Value = 15
$ws.Range("F536").Value = 100112023
$ws.Range("G536").Value = "Brócoli"
$ws.Range("H536").Value = "Sin especificar"
$ws.Range("I536").Value = "Segunda"
$ws.Range("J536").Value = 1200
$ws.Range("K536").Value = 500
$ws.Range("L536").Value = 600
$ws.Range("M536").Value = 550
$ws.Range("N536").Value = "`$/unidad"
$ws.Range("O536").Value = "Región de Arica y Parinacota"
$ws.Range("P536").Value = 550
$ws.Range("Q536").Value = 1
$ws.Range("R536").Value = "Hortaliza"
